$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 159; this shifts the existing rows 159-191
# down to become rows 160-192, matching the diff's row-shift pattern.
$ws.Rows("159:159").Insert()

# Populate the newly inserted row 159 with the new weekly record.
$ws.Range("A159").Value = 11
$ws.Range("B159").Value = "Vega Monumental Concepción"
$ws.Range("C159").Value = "Bíobío"
$ws.Range("D159").Value = 44995
$ws.Range("E159").Value = 8
$ws.Range("F159").Value = 100112043
$ws.Range("G159").Value = "Pepino ensalada"
$ws.Range("H159").Value = "Sin especificar"
$ws.Range("I159").Value = "Primera"
$ws.Range("J159").Value = 100
$ws.Range("K159").Value = 6000
$ws.Range("L159").Value = 6500
$ws.Range("M159").Value = 6250
$ws.Range("N159").Value = "$/caja 60 unidades"
$ws.Range("O159").Value = "Provincia de Limarí"
$ws.Range("P159").Value = 104
$ws.Range("Q159").Value = 60
$ws.Range("R159").Value = "Hortaliza"
